$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "27.558.48"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.10%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.755.04"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.30%  "

$ws.Cells.Item(4, 5).Value = "  -0.15%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "324.95"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.21%  "

$ws.Cells.Item(6, 5).Value = "  -0.15%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.4581"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +2.35%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.3562"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.51%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.07455"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.37%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "41.50"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.76%  "

$ws.Cells.Item(11, 5).Value = "  -1.24%  "

$ws.Cells.Item(12, 5).Value = "  -0.16%  "

$ws.Cells.Item(13, 5).Value = "  +0.73%  "

$ws.Cells.Item(14, 5).Value = "  -0.47%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "7.163"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.24%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "1.755.34"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.20%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "93.50"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.01%  "

$ws.Cells.Item(18, 5).Value = "  -0.74%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.06419"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.51%  "

$ws.Cells.Item(21, 5).Value = "  +1.23%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "5.741"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -2.08%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "27.604.73"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.05%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "11.22"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.19%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.071"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -1.17%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "165.14"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +2.05%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "20.12"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.41%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "1.955.26"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.16%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "2.134"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +1.13%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "125.30"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.14%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "1.078"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.39%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.09222"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +2.32%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "3.664"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.17%  "

$ws.Cells.Item(34, 5).Value = "  -0.49%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "11.74"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -2.03%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.02275"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -1.89%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.06018"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.59%  "

$ws.Cells.Item(38, 5).Value = "  -0.26%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.6266"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.27%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "4.924"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.47%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "1.182"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -2.09%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "1.383"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.74%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "7.760"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.07%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "13.22"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.32%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "3.719"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.12%  "

$ws.Cells.Item(46, 5).Value = "  -0.30%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "121.84"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.38%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.936"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.57%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.06896"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.34%  "

$ws.Cells.Item(50, 5).Value = "  -2.30%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "72.05"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.29%  "
